$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 63 - this shifts existing rows 63:164 down to 64:165
$ws.Rows("63:63").Insert()

# Populate the newly inserted row 63 with the new record's data.
# Columns A,B,C,E,F,G,Q,R are constant across all data rows in this sheet.
$ws.Range("A63").Value = 11
$ws.Range("B63").Value = "Vega Monumental Concepción"
$ws.Range("C63").Value = "Bíobío"
$ws.Range("D63").Value = 44580
$ws.Range("E63").Value = 8
$ws.Range("F63").Value = 100114001
$ws.Range("G63").Value = "Papa"
$ws.Range("H63").Value = "Patagonia"
$ws.Range("I63").Value = "1a nueva(o)"
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 8500
$ws.Range("L63").Value = 9000
$ws.Range("M63").Value = 8750
$ws.Range("N63").Value = "$/saco 25 kilos"
$ws.Range("O63").Value = "Región de La Araucanía"
$ws.Range("P63").Value = 350
$ws.Range("Q63").Value = 25
$ws.Range("R63").Value = "Hortaliza"
